$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.253.79'
$ws.Range("E2").Value = '  +0.49%  '

$ws.Range("D3").Value = '3.131.20'
$ws.Range("E3").Value = '  +0.65%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.77'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.12%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.37'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.66%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.523'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.22%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.156'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.09%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.41'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.50%  '

$ws.Range("E11").Value = '  -0.31%  '

$ws.Range("E12").Value = '  +0.55%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.46'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.18%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.122'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.56%  '

$ws.Range("D15").Value = '3.650.36'
$ws.Range("E15").Value = '  +0.67%  '

$ws.Range("D16").Value = '67.214.35'
$ws.Range("E16").Value = '  +0.48%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.15'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.83%  '

$ws.Range("D18").Value = '3.130.67'
$ws.Range("E18").Value = '  +0.67%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.39'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.54%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '494.09'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.80%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.710'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.86%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.91'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.99%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '84.30'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.14%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.34'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.30'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.20%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.44'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.81%  '

$ws.Range("E27").Value = '  -0.01%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.95'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.41%  '

$ws.Range("E29").Value = '  -2.13%  '

$ws.Range("E30").Value = '  +0.25%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '28.80'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.69%  '

$ws.Range("E32").Value = '  -0.74%  '

$ws.Range("D33").Value = '0.0₃0948'
$ws.Range("E33").Value = '  -6.27%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.03%  '

$ws.Range("E35").Value = '  +0.06%  '

$ws.Range("E36").Value = '  -2.38%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '46.98'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.91%  '

$ws.Range("E38").Value = '  -2.50%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '50.15'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.14%  '

$ws.Range("E40").Value = '  -1.21%  '

$ws.Range("E41").Value = '  +1.71%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.58'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.18%  '

$ws.Range("D43").Value = '2.837.23'
$ws.Range("E43").Value = '  -0.02%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '387.14'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.33%  '

$ws.Range("E45").Value = '  -6.82%  '

$ws.Range("E46").Value = '  -1.98%  '

$ws.Range("E47").Value = '  +0.30%  '

$ws.Range("E48").Value = '  -0.03%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.06'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.03%  '

$ws.Range("E50").Value = '  -0.16%  '

$ws.Range("E51").Value = '  -0.23%  '
